$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: new departure entry (Sunday, Jan 15 - Wizz Air A321 G-WUKI)
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "Sunday, Jan 15"
$ws.Range("C25").Value = "10:10 AM"
$ws.Range("D25").Value = "W92065"
$ws.Range("E25").Value = "London"
$ws.Range("F25").Value = "(LTN)"
$ws.Range("G25").Value = "Wizz Air "
$ws.Range("H25").Value = "A321"
$ws.Range("I25").Value = "(G-WUKI)"
$ws.Range("J25").Value = "10:17 AM"
$ws.Range("K25").Font.Size = 11
$ws.Range("L25").Value = "0 hours, 7 minutes"
$ws.Range("M25").Font.Size = 11

# Row 26: new departure entry (Sunday, Jan 15 - Wizz Air A321 HA-LTC)
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "Sunday, Jan 15"
$ws.Range("C26").Value = "12:40 PM"
$ws.Range("D26").Value = "W62090"
$ws.Range("E26").Value = "Eindhoven"
$ws.Range("F26").Value = "(EIN)"
$ws.Range("G26").Value = "Wizz Air "
$ws.Range("H26").Value = "A321"
$ws.Range("I26").Value = "(HA-LTC)"
$ws.Range("J26").Value = "12:42 PM"
$ws.Range("K26").Font.Size = 11
$ws.Range("L26").Value = "0 hours, 2 minutes"
$ws.Range("M26").Font.Size = 11

Write-Output "rows added"
